# Auto-generated update of Tiamat_Profits leve-profit computed columns (H-N)
# per scheduled-runner market-price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 666.6799999999999
$ws.Range("J41").Value = 552.53845
$ws.Range("L41").Value = 552.53845
$ws.Range("N41").Value = -1432.53845

$ws.Range("H69").Value = 3678349.2
$ws.Range("I69").Value = 4903794
$ws.Range("J69").Value = 2015
$ws.Range("K69").Value = 14711382
$ws.Range("L69").Value = 6045
$ws.Range("M69").Value = -14710508
$ws.Range("N69").Value = -7793

$ws.Range("H72").Value = 3678349.2
$ws.Range("I72").Value = 4903794
$ws.Range("J72").Value = 2015
$ws.Range("K72").Value = 44134146
$ws.Range("L72").Value = 18135
$ws.Range("M72").Value = -44129778
$ws.Range("N72").Value = -26871

$ws.Range("H74").Value = 1964399.9
$ws.Range("I74").Value = 2503521.5
$ws.Range("K74").Value = 2503521.5
$ws.Range("M74").Value = -2502585.5

$ws.Range("H77").Value = 1964399.9
$ws.Range("I77").Value = 2503521.5
$ws.Range("K77").Value = 12517607.5
$ws.Range("M77").Value = -12512927.5

$ws.Range("H80").Value = 3955144.5
$ws.Range("I80").Value = 2278.2
$ws.Range("J80").Value = 6995811
$ws.Range("K80").Value = 6834.599999999999
$ws.Range("L80").Value = 20987433
$ws.Range("M80").Value = -5836.599999999999
$ws.Range("N80").Value = -20989429

$ws.Range("H83").Value = 3955144.5
$ws.Range("I83").Value = 2278.2
$ws.Range("J83").Value = 6995811
$ws.Range("K83").Value = 20503.8
$ws.Range("L83").Value = 62962299
$ws.Range("M83").Value = -15511.8
$ws.Range("N83").Value = -62972283

$ws.Range("H137").Value = 19755.547
$ws.Range("I137").Value = 22101.299
$ws.Range("J137").Value = 1380.5
$ws.Range("K137").Value = 66303.897
$ws.Range("L137").Value = 4141.5
$ws.Range("M137").Value = -63753.897
$ws.Range("N137").Value = -9241.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1637515.9
$ws.Range("I32").Value = 1868158.1
$ws.Range("J32").Value = 125528.445
$ws.Range("K32").Value = 1868158.1
$ws.Range("L32").Value = 125528.445
$ws.Range("M32").Value = -1867871.1
$ws.Range("N32").Value = -126102.445

$ws.Range("H76").Value = 49525.668
$ws.Range("J76").Value = 49525.668
$ws.Range("L76").Value = 49525.668
$ws.Range("N76").Value = -50201.668

$ws.Range("H79").Value = 49525.668
$ws.Range("J79").Value = 49525.668
$ws.Range("L79").Value = 49525.668
$ws.Range("N79").Value = -51865.668

$ws.Range("H102").Value = 1816.6666
$ws.Range("I102").Value = 1863.6364
$ws.Range("J102").Value = 1300
$ws.Range("K102").Value = 1863.6364
$ws.Range("L102").Value = 1300
$ws.Range("M102").Value = -241.6364000000001
$ws.Range("N102").Value = -4544

$ws.Range("H122").Value = 806.4211
$ws.Range("I122").Value = 521.8333
$ws.Range("J122").Value = 1294.2858
$ws.Range("K122").Value = 1565.4999
$ws.Range("L122").Value = 3882.8574
$ws.Range("M122").Value = 884.5001
$ws.Range("N122").Value = -8782.857400000001

$ws.Range("H132").Value = 2548876.2
$ws.Range("I132").Value = 2685141.8
$ws.Range("K132").Value = 8055425.399999999
$ws.Range("M132").Value = -8052895.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 14207.333
$ws.Range("J88").Value = 10000
$ws.Range("L88").Value = 10000
$ws.Range("N88").Value = -10812

$ws.Range("H91").Value = 14207.333
$ws.Range("J91").Value = 10000
$ws.Range("L91").Value = 10000
$ws.Range("N91").Value = -12808

$ws.Range("H99").Value = 566.0769
$ws.Range("I99").Value = 588
$ws.Range("J99").Value = 524.6667
$ws.Range("K99").Value = 588
$ws.Range("L99").Value = 524.6667
$ws.Range("M99").Value = 910
$ws.Range("N99").Value = -3520.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34193.977
$ws.Range("I31").Value = 45340.89
$ws.Range("J31").Value = 12696.357
$ws.Range("K31").Value = 45340.89
$ws.Range("L31").Value = 12696.357
$ws.Range("M31").Value = -45045.89
$ws.Range("N31").Value = -13286.357

$ws.Range("H34").Value = 34193.977
$ws.Range("I34").Value = 45340.89
$ws.Range("J34").Value = 12696.357
$ws.Range("K34").Value = 45340.89
$ws.Range("L34").Value = 12696.357
$ws.Range("M34").Value = -45138.89
$ws.Range("N34").Value = -13100.357

$ws.Range("H132").Value = 2619.2173
$ws.Range("I132").Value = 1940
$ws.Range("J132").Value = 3502.2
$ws.Range("K132").Value = 5820
$ws.Range("L132").Value = 10506.6
$ws.Range("M132").Value = -3290
$ws.Range("N132").Value = -15566.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 53.26923
$ws.Range("I12").Value = 105.5
$ws.Range("J12").Value = 37.6
$ws.Range("K12").Value = 316.5
$ws.Range("L12").Value = 112.8
$ws.Range("M12").Value = -143.5
$ws.Range("N12").Value = -458.8

$ws.Range("H68").Value = 12041.333
$ws.Range("I68").Value = 33834
$ws.Range("J68").Value = 1145
$ws.Range("K68").Value = 101502
$ws.Range("L68").Value = 3435
$ws.Range("M68").Value = -100691
$ws.Range("N68").Value = -5057

$ws.Range("H71").Value = 12041.333
$ws.Range("I71").Value = 33834
$ws.Range("J71").Value = 1145
$ws.Range("K71").Value = 304506
$ws.Range("L71").Value = 10305
$ws.Range("M71").Value = -300450
$ws.Range("N71").Value = -18417

$ws.Range("H86").Value = 600
$ws.Range("J86").Value = 600
$ws.Range("L86").Value = 1800
$ws.Range("N86").Value = -4172

$ws.Range("H89").Value = 600
$ws.Range("J89").Value = 600
$ws.Range("L89").Value = 5400
$ws.Range("N89").Value = -17256

$ws.Range("H132").Value = 1995
$ws.Range("J132").Value = 1995
$ws.Range("L132").Value = 17955
$ws.Range("N132").Value = -23015

$ws.Range("H141").Value = 2960.8
$ws.Range("I141").Value = 888.75
$ws.Range("J141").Value = 6644.4443
$ws.Range("K141").Value = 2666.25
$ws.Range("L141").Value = 19933.3329
$ws.Range("M141").Value = 2513.75
$ws.Range("N141").Value = -30293.3329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 47586.316
$ws.Range("I132").Value = 1390.6
$ws.Range("J132").Value = 146577.14
$ws.Range("K132").Value = 4171.799999999999
$ws.Range("L132").Value = 439731.42
$ws.Range("M132").Value = -1641.799999999999
$ws.Range("N132").Value = -444791.42

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 829
$ws.Range("I16").Value = 702.1539
$ws.Range("J16").Value = 1158.8
$ws.Range("K16").Value = 702.1539
$ws.Range("L16").Value = 1158.8
$ws.Range("M16").Value = -532.1539
$ws.Range("N16").Value = -1498.8

$ws.Range("H132").Value = 310551.97
$ws.Range("J132").Value = 591467.25
$ws.Range("L132").Value = 1774401.75
$ws.Range("N132").Value = -1779461.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 852.0714
$ws.Range("I100").Value = 1009.3333
$ws.Range("J100").Value = 809.1818
$ws.Range("K100").Value = 2018.6666
$ws.Range("L100").Value = 1618.3636
$ws.Range("M100").Value = -1477.6666
$ws.Range("N100").Value = -2700.3636

$ws.Range("H132").Value = 10717.583
$ws.Range("I132").Value = 2601.3333
$ws.Range("J132").Value = 18833.834
$ws.Range("K132").Value = 7803.999899999999
$ws.Range("L132").Value = 56501.50199999999
$ws.Range("M132").Value = -5273.999899999999
$ws.Range("N132").Value = -61561.50199999999
